# TD correction in test_LoginSkillBoard.py:
# Column A ("TestCase") values get corrected from the generic "Sign_Up"
# placeholder to per-row test-case ids (tc_1, tc_3, tc_4, tc_5, tc_6) and,
# for the row that used to hold test_2's data, the real case name
# "test_signup_on_skillboard". All other columns (username/password/
# firstname/lastname/email) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NB: write row 3 ("test_signup_on_skillboard") last so the workbook's
# shared-string table regenerates new entries in the same order as the
# target file (tc_1, tc_3, tc_4, tc_5, tc_6, test_signup_on_skillboard).
$ws.Range("A2").Value = "tc_1"
$ws.Range("A4").Value = "tc_3"
$ws.Range("A5").Value = "tc_4"
$ws.Range("A6").Value = "tc_5"
$ws.Range("A7").Value = "tc_6"
$ws.Range("A3").Value = "test_signup_on_skillboard"

# The longer text in column A (and the now-unused "Sign_Up" string that
# used to govern its width) means the column needs to be resized - and
# with it, Excel recomputes the best-fit width for the sheet's other
# text columns too.
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(6).AutoFit()

# The data column (F) no longer carries the extra "applyFill" cell style -
# drop back to the workbook's default/Normal style.
$ws.Range("F2:F7").Style = "Normal"

# Cursor ends up one row higher than before.
$null = $ws.Range("D10").Select()
